$wb = $excel.ActiveWorkbook

# --- Step 1: structural sheet changes first. Worksheet handles in this
# runtime are position-based, so a handle fetched before a sheet is
# inserted/moved can silently end up pointing at the wrong sheet afterwards.
# Always grab a fresh handle (by name) immediately before it's needed. ---

$wsForecast = $wb.Worksheets.Add()
$wsForecast.Name = "PO Forecast"

# Fetch "Monthly Trend" only now (after Add shifted indices) so Move lands
# the new sheet right after it, giving the final tab order:
#   Weekly Quantity, Monthly Trend, PO Forecast
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsForecast.Move($null, $wsMonthly)

# --- Step 2: re-fetch fresh references by name now that sheet order/identity
# is final, and make the cell edits ---

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

$wsForecast = $wb.Worksheets.Item("PO Forecast")

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Data rows
$data = @(
    @(45326.99999999999, 44, 30.05041756171731, 58.46372765955587),
    @(45333.99999999999, 42, 28.46976486618497, 56.96038565479374),
    @(45340.99999999999, 40, 25.4486244909025, 54.18239794162753),
    @(45361.99999999999, 34, 20.30049208752947, 47.1463078639993),
    @(45368.99999999999, 31, 16.76993657043973, 45.76321037110065),
    @(45375.99999999999, 29, 15.05890782814753, 44.76141679845775),
    @(45382.99999999999, 27, 12.18082910701204, 41.67516673572469),
    @(45389.99999999999, 25, 11.34446765966521, 40.64777671736104),
    @(45396.99999999999, 23, 8.262574460752331, 36.87559519684725),
    @(45403.99999999999, 21, 6.504319690937464, 34.71526951356827),
    @(45410.99999999999, 19, 3.999254132605749, 32.65517293342327),
    @(45417.99999999999, 16, 1.177958577338192, 31.44157650724085)
)

$row = 2
foreach ($r in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $r[0]
    $wsForecast.Cells.Item($row, 2).Value = $r[1]
    $wsForecast.Cells.Item($row, 3).Value = $r[2]
    $wsForecast.Cells.Item($row, 4).Value = $r[3]
    $row++
}

# --- Step 3: formatting, matching the look of the existing sheets ---
# Header row: bold / centered / thin-bordered, same as the other sheets'
# header row (copy format from an existing header cell so the same style
# index is reused rather than a new one being created).
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)   # xlPasteFormats

# Date column (A2:A13): same date number format as the other sheets' date
# column.
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A13").PasteSpecial(-4122)  # xlPasteFormats
